# loss function changed to 0.7xMSE+0.3xIoU
# Adds a new "scheduler change" row (18) documenting commit 1dcfe88, updates the
# "current best commit" pointer (M2) to that same commit, and starts a new row
# (19) for the in-progress loss-function change (0.7*MSE + 0.3*IoU).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 18: scheduler tweak (commit 1dcfe88) -----------------------------
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = 40
$ws.Range("D18").Value = 12
$ws.Range("F18").Value = "Изменён scheduler для снижения скорости уменьшения learning rate"
$ws.Range("G18").Value = "параметры теста 4"
$ws.Range("I18").Value = "1dcfe88"
$ws.Range("H18").Value = "Train IoU: 0.56, Val IoU: 0.55. Изменения незначительные, но пойдёт. Дальнейшее послабление оптимизатора приведёт к переобучению."

# Pointer to the current best commit now follows row 18's result.
$ws.Range("M2").Value = "1dcfe88" + "`n"

# --- Row 19: new loss-function experiment (in progress) -------------------
$ws.Range("B19").Value = 1
$ws.Range("C19").Value = 40
$ws.Range("D19").Value = 13
$ws.Range("F19").Value = "Изменение функции потерь на комбинацию MSE и IoU"
$ws.Range("G19").Value = "параметры теста 4"

# Row heights to accommodate the wrapped text added above.
$ws.Rows.Item(18).RowHeight = 75
$ws.Rows.Item(19).RowHeight = 30

# Reflect where the user was last looking/working in the sheet.
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("G18").Select()
